$wb = $excel.ActiveWorkbook

# --- Department sheet: remove the two blank rows above the table ---
$wsDept = $wb.Worksheets.Item("Department")
$wsDept.Rows("1:2").Delete()

# --- Grade sheet: remove the two blank rows above the table ---
$wsGrade = $wb.Worksheets.Item("Grade")
$wsGrade.Rows("1:2").Delete()

# --- Restore/update selections on each sheet ---
$wsEmployee = $wb.Worksheets.Item("Employee")
$wsEmployee.Activate()
$wsEmployee.Range("F7").Select()

$wsDept.Activate()
$wsDept.Rows("1:2").Select()

$wsGrade.Activate()
$wsGrade.Range("F11").Select()
